# Daily update of covid19 tracker data files
# Targeted shared-string text edits on the 'Country Updates' sheet,
# applied via Range.Value assignments (cell addresses derived from the
# original workbook's shared-string usage).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

$text_F5 = @'
All shops and public spaces are closed from March 20 until 12 April, only supermarkets, small groceries shops, hardware stores and pharmacies as well as certain parts of public administration remain open. As of 3 April, the following additional activities are excluded from the mandatory quarantine: distribution of construction materials, mining, forrestry, manufacturing of wood, furniture, leather as well as agricultural or construction machines, imports and exports of already produced items, nuclear fuel production, maintenance and desinfection services as well as a minimum service of banks and savings and credit cooperatives.
'@
$ws.Range("F5").Value = $text_F5

$text_F6 = @'
The government (National Cabinet, made up of the Prime Minister, premiers and territory leaders) has now limited both indoor and outdoor gatherings to two persons only. Exceptions to this limit include: i) People of the same household going out together; ii) Funerals (a maximum of 10 people); iii) Wedding (a maximum of 5 people); iv) Family units. The government’s strong guidance to all Australians is to stay home unless for: i) shopping for what you need (food and necessary supplies); ii) medical or health care needs, including compassionate requirements; iii) exercise in compliance with the public gathering requirements; iv) work and study if you can’t work or learn remotely. However, each State Government has adopted rules that are unique to their situation. In all jurisdictions, playgrounds, skate parks and outside gyms in public places have been closed.
'@
$ws.Range("F6").Value = $text_F6

$text_J6 = @'
New Commonwealth Government measures include support for business investment, including by increasing the threshold for assets eligible for instant tax write-off and expanding access to this regime to larger firms (with annual turnover of <500 million AUD up from <50 million AUD), as well as through new accelerated depreciation measures. Measures have also been implemented to support business cash flows (including through payments to SMEs that employ workers to help finance payments of salaries and wages plus wage subsidies to SMEs to support them in retaining apprentices and trainees). 
A range of measures have been taken to allow firms that are affected by the pandemic to defer payment of tax liabilities for up to 4 months without incurring interests or penalties for late payment (incl. for income taxes, goods and services tax (GST), excises and fringe benefits tax). These have been complemented with measures to accelerate tax refunds, notably of excess input-GST credits.
Under a new Coronavirus SME Guarantee Scheme, the Government will guarantee 50 per cent of new loans issued by eligible lenders to SMEs. The total lending capacity of the facility will be A$40bn (2% of GDP). Under a plan put forward by the banking industry, businesses with up to A$10mn in total loan facilities will be able to defer their loan repayments for six months.

'@
$ws.Range("J6").Value = $text_J6

$text_C20 = @'
• Mandated nation-wide quarantine for at least 15 days since March 15th, which has been further extended until April 11th. 
• Suspension of all non-essential work from March 30 to April 9, inclusive.
'@
$ws.Range("C20").Value = $text_C20

$text_I20 = @'
March 30: Workers affected by the suspension of all non-essential activities  will receive their full salary during the suspension, and will be required to compensate the work hours before December 31, 2020.

March 27: Dismissals for reasons related to Covid19 will not be considered justified rom March 27 until the end of the health crisis. Temporary contracts can not be cancelled (they can be interrupted during the crisis but must be resumed after the end of the lockdown). 

March 17:
• Changes in the temporary employment adjustment schemes (ERTEs) in order to avoid outright dismissal by promoting temporary unemployment (the suspension of employment) or reductions in working time. In particular: (i) the temporary employment adjustment schemes have been significantly simplified, (ii) access conditions to these schemes have been extended, (iii) unemployment benefits received under the temporary employment adjustment scheme do not count in terms of consumption of unemployment benefit rights, and (iv) as further detailed under firm-specific measures, employers will be exempt from social contributions during the period of application of the temporary employment adjustment scheme. 
• An extraordinary allowance is provided for self-employed workers affected by the suspension of economic activity (70% of the SSC base, at least for one month). 
• Additional budgetary funds of EUR 300 million to ensure the provision of assistance to dependent persons.  
• Additional flexibility for local authorities to use their 2019 budgetary surplus to fund social services and primary assistance to dependent persons (EUR 300 million) 
• The social benefit for energy provision (‘bono social’) will be automatically extended until September 15.

• March 12: Supplemental credit of EUR 25 million to cover meal allowances to ensure the basic access to food for vulnerable children affected by the suspension of educational activity in schools.                                                                                                                                                                                                         
• March 10: Increased sick pay for coronavirus infected workers or those quarantined, from 60 to 75 percent of the regulatory base, paid by the social security budget.
Other employment and social measures:
• Three-month credit postponement on mortgage payments for the most vulnerable. (March 17, expanded on March 31)
• Broadened scope for protected families in the supply of water and energy. Telecommunication services are also guaranteed. (March 17, expanded on March 31)
• Increased flexibility of working conditions, encouragement of telework and adjustment of working times. Reduced working times are permitted for workers having to take care of children, elderly or dependent persons. (March 17)

'@
$ws.Range("I20").Value = $text_I20

$text_J20 = @'
March 17:
• Exemption of social security contributions, directed to maintain employment in temporary employment adjustment schemes (ERTEs) due to COVID -19: 100% for SMEs, 75% of employer’s social security contributions for other companies. 
• The Government will open a new line of guarantees via the national development bank (ICO) of up to EUR 100 billion so that the financial sector provides liquidity to firms and self-employed workers to fund working capital, payment of bills and other needs to maintain operations and protect economic activity and employment. The main features of the first tranche of guarantees amounting to EUR 20 billion are the following: (i) 50% of the tranche will cover loans to SMEs and self-employed workers, and (ii) the guarantee will cover 80% of new loans and loan renewals to SMEs and self-employed workers, whereas for the rest of companies, the guarantee will cover 70% of new loans and 60% of loan renewals.
• The ICO is allowed to increase its funding by EUR 10 billion to extend its existing lines of credit to companies and self-employed workers. 
• Additional guarantees of up to EUR 2 billion through the Spanish Export Insurance Credit Company (CESCE), for financial institutions to provide new working capital credit to companies.
• Guarantees for loan maturity extensions to farmers using credit lines under the special 2017 drought credit line. 
• Support to the digitalisation of small and medium companies through grants and loans to finance investment in digital equipment or solutions for remote working conditions (programme ACELERA PYME). 

March12:
• Increased flexibility of the tax deferral regime: possibility for SMEs and self-employed workers, upon request, to defer tax payments for six months, and benefit from interest rate subsidies. This measure will allow for an injection of up to EUR 14 billion in the system. 
• 50% exemption from employers social security contributions, from February to June 2020, for workers with permanent discontinuous contracts in the tourism sector and related activities. 
• Specific ICO (Instituto de Crédito Oficial) financing facility amounting to EUR 400 million to support, through liquidity provision, firms and self-employed workers in the tourism sector affected by COVID-19. 
• Deferral of the reimbursement of principal and/or interest of loans received from the Ministry of Industry, Trade and Tourism if COVID-19 has caused disruptions to the beneficiary firm. 

Other firm-specific measures:
• Prior government authorisation for third country FDI in strategic sectors.
• Acceleration of custom procedures for imports and exports in the industrial sector for six months.
• Facilitation of procedural and administrative requirements for companies. (March 17) 
• Measures to support the aviation sector: EU Commission has been requested to provide urgently flexibility for airport slots. 
• Publication of sectoral guidelines on how to operate in labour related aspects in the context of coronavirus, and on good practices for businesses and workers in the tourism sector.
• Introduction of a special regime for the suspension of public contracts, foreseeing compensation in certain cases. (March 17)



'@
$ws.Range("J20").Value = $text_J20

$text_G23 = @'
• On 20 March, the government announced about EUR 39 million (0.02% of GDP) of extra spending on health care to cover expenses related to the coronavirus. The government also plans additional funding to cover local governments’ spending. A raft of emergency measures was adopted to safeguard medicines and medical supplies. The sale of medical supplies and pharmaceuticals used in the public health-care system is restricted and municipalities are no longer obliged to provide timely non-emergency care (until 13 May). 

'@
$ws.Range("G23").Value = $text_G23

$text_D39 = @'
- To this day, no travel bans have been adopted by the Mexican government.
- Some airlines, such as Aeromexico imposed 14-day quarantine to passengers arriving from Europe and China, before being able to fly to the United States. Other airlines recommend passengers arriving from high-contagious-rate countries to voluntary quarantine for a two-week period.
- Mexico and the United States agreed to restrict non-essential border crosses (including for tourism and shopping purposes) through their common border, in order to keep bilateral trade flow during the sanitary crisis. Besides commercial exchanges, workers authorized to work on both sides of the border, as well as medical and emergency services are exempted from the aforementioned restrictions.
- On 17 March, Guatemala closed its border with Mexico for an initial period of 15 days and Belize limited all border crossing with Mexico on 21 March, except for the freight-transport sector.

'@
$ws.Range("D39").Value = $text_D39

$text_G39 = @'
- On 24 March, the Federal Government announced a hospital conversion to better coordinate and integrate available space in hospitals and other infrastructures to treat patients with COVID-19. This will include the use of the military health system, as well as the Plan DN-III and Plan Marina for catastrophic events. 
- The estimated human resources for the three stages of hospital conversion are the following:
First stage: 1 025 available beds, 4 337 nurses and 1 811 doctors. 
Second stage: the number of available beds that can adapt a ventilator is 5 059, with 18 873 nurses and 6 845 doctors. 
Third stage: the number of extra beds available is 2 342, with additional 9 908 nurses and 3 982 doctors. 
- On 1st of April, the Ministry of Health announced the launch of the app called “COVID-19MX”, that seeks to reinforce knowledge and promote health on the matter of COVID-19, as well as to facilitate exchange of information between health officials and citizens. The app gathers data on the health status of its users and determines whether they belong to risk groups. Additionally, users can use this app for self-diagnosis purposes.

'@
$ws.Range("G39").Value = $text_G39

$text_H39 = @'
- On 24 March, the extraordinary budget was allocated to the Ministry of National Defence (up to MXN 4 billion/USD 160 million) and to the Ministry of the Navy (up to MXN 500 million/USD 20.5 million), to reinforce the Plan DN-III and Plan Marina in response to the COVID emergency.
- The Health Institute for Wellbeing (Instituto de Salud para el Bienestar) authorized MXN 4.5 billion to buy medical devices, medicine and health equipment.
- Through the Health and the Wellbeing Institute, States will receive transfers up to MXN 15.3 billion (USD 629 million). In addition, States will have advance funding for MXN 10 billion (USD 411 million).   
- Various States have been taken fiscal measures such as deferral of tax payments, discounts in payroll taxes and property taxes, or reduction in special taxes like vehicular ownership, etc. as well as postponements of tax inspection.
'@
$ws.Range("H39").Value = $text_H39

$text_J39 = @'
- The Federal Government, through the development bank Nacional Financiera (Nafin), will support SMEs up to MXN 36.3 billion (USD 1 billion). One million credit products of MXN 25 000 each (USD 1 000); and 500 000 credit products will be provided to the formal economy and 500 000 loans will be extended to the informal economy. The loans will be settled over a period of 3 years, with no payment in the first three months and further monthly payments of MXN 1 000 (USD 42), at an average rate of 6.5% per year. The Government will further boost the economy by providing new financing of at least MXN 24.8 billion via other development banks. 
- Some States are taking different actions to support SMEs in the formal and informal sectors: Discounts on payroll taxes for SME’s and credit lines for payment of payrolls or supply merchandise for SMEs, among others. Many credits will be settled over a period of 2 years, with no payment during the first three or four months and low interest rates. Financial support to self-employed, rural workers (corn producers), family businesses and vulnerable groups (women, indigenous people) has also been put in place.

'@
$ws.Range("J39").Value = $text_J39

$text_K39 = @'
- On 20 March, Banco de México decided to bring forward its scheduled March monetary policy decision and lowered the overnight interbank interest rate target by 50 basis points to 6.5%.
- In collaboration with Ministry of Finance, Banco de México is strengthening the Market Makers Program for government securities. Banco de México, acting as the Federal Government’s financial agent, will 1) allow market makers to unwind their short positions in Banco de México’s government securities lending facility with other government securities available in the market and 2) extend by one day the exercise period of the green-shoe option. 
- Banco de México and the US Federal Reserve established a temporary US dollar swap line of USD 60 billion on 19 March. This swap line will fund US dollar auctions carried out by Banco de México among credit institutions to increase the availability of USD funding for the private sector.
- As a response to worsening market conditions in the Mexican peso, the Foreign Exchange Commission announced on 9 March a USD 10 billion increase in the total amount available of the Non-Deliverable Forward (NDF) program, from USD 20 billion to USD 30 billion.
- Banco de México reduced by MXN 50 billion the amount of the Monetary Regulation Deposit (DRM) to support active transactions of commercial and development banks, thus improving their liquidity and their capacity to grant credits and to maintain or expand their credit lines.
- To strengthen commercial banks’ financing alternatives and to ensure the necessary liquidity to perform their operations adequately, Banco de México adjusted the Ordinary Additional Liquidity Facility interest rate from a factor of 2-2.2 to 1.1 times Banco de México’s target for the overnight interbank interest rate. 
- The National Banking and Securities Commission (CNBV) published temporary accounting criteria for banks, credit unions and institutions related to popular finance so as to facilitate the maintenance of credit lines. 
- In addition, CNBV issued a recommendation letter to banks to avoid share-buybacks or any other compensation measure for shareholders, coupled with dividend distribution restriction for 2019 and 2020 tax years. The recommendation aims to redirect resources and strengthen the institutions to be able to support the economy.

'@
$ws.Range("K39").Value = $text_K39

$text_G51 = @'
- The central government will compensate municipalities and regions for extraordinary measures and extra costs in health and medical care associated with the COVID-19 virus. On 2 April, the amount was raised to SEK 3bn from SEK 2bn. 
- On 2 April, the government announced that grants to municipalities and regions will increase by SEK 20bn in 2020, up from SEK 5bn previously planned, of which SEK 12.5bn will be permanent. 
- Sickness benefit to be paid for the first day of sickness.
'@
$ws.Range("G51").Value = $text_G51

$text_H51 = @'
The government has presented three comprehensive packages amounting to SEK 661.8 billion (about EUR 62.1 billion, nearly 13% of GDP), including liquidity provisions and credit guarantees if these are used by firms to the maximum.
'@
$ws.Range("H51").Value = $text_H51

